# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to reflect the latest generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row number -> new F-column value }
$sheetUpdates = @{
    "展览" = @{
        2  = 1317
        4  = 14453
        5  = 17042
        8  = 45
        22 = 57
        23 = 7
        24 = 7014
        28 = 22
        30 = 40
        31 = 5812
        32 = 127
        34 = 211
        35 = 4966
    }
    "全部类型" = @{
        2  = 1317
        4  = 14453
        5  = 17042
        8  = 45
        23 = 57
        24 = 7
        25 = 7014
        29 = 22
        31 = 40
        33 = 5812
        34 = 127
        36 = 211
        37 = 4966
    }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetUpdates[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $ws.Cells.Item($rowNum, 6).Value = $rows[$rowNum]
    }
}
